# R22 UAT 1 data refresh for stocksReceivedData + selection/active-tab cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stocksReceivedData")

# --- New serial / invoice numbers for rows 2-5 ----------------------------
$ws.Range("B2").Value = 21585
$ws.Range("C2").Value = 21590
$ws.Range("E2").Value = 5465627

$ws.Range("B3").Value = 21591
$ws.Range("C3").Value = 21595
$ws.Range("E3").Value = 5465628

$ws.Range("B4").Value = 21596
$ws.Range("C4").Value = 21600
$ws.Range("E4").Value = 5465629

$ws.Range("B5").Value = 21601
$ws.Range("C5").Value = 21605
$ws.Range("E5").Value = 5465630

# Rows 2-4 now get the same "thick bottom" row look row 5 already has.
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75

# Amount column (E) grows a medium/thick bottom border on rows 2-4, matching
# the border already used on row 5's cells.
foreach ($r in 2, 3, 4) {
    $rng = $ws.Range("E" + $r)
    $rng.Borders.Item(7).LineStyle = 1   # left
    $rng.Borders.Item(7).Weight = 2
    $rng.Borders.Item(8).LineStyle = 1   # top
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(10).LineStyle = 1  # right
    $rng.Borders.Item(10).Weight = 2
    $rng.Borders.Item(9).LineStyle = 1   # bottom
    $rng.Borders.Item(9).Weight = -4138
}

# Serial NumS/NumE columns (B/C) move onto the plain thin border used
# throughout the rest of the table.
foreach ($ref in "B3", "B4", "B5", "C5") {
    $rng = $ws.Range($ref)
    $rng.Borders.Item(7).LineStyle = 1   # left
    $rng.Borders.Item(7).Weight = 2
    $rng.Borders.Item(8).LineStyle = 1   # top
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(9).LineStyle = 1   # bottom
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(10).LineStyle = 1  # right
    $rng.Borders.Item(10).Weight = 2
}

# --- Selection / active sheet bookkeeping ---------------------------------
# Previously the "stocksEnquiryWorkingStock" tab was left active; this UAT
# pass leaves "stocksReceivedData" selected at C5 instead.
$ws.Activate()
$ws.Range("C5").Select()

Write-Host "StockManagement UAT1 refresh applied"
